$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.597.51"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.60"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.47"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4645"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3880"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07877"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9719"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.94"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.897.60"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.696"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.962"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06991"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.12"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001002"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.81"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.599.11"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.295"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.98"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.113"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.084.65"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.73"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.19"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.717"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.981"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.14"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09328"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9154"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.259"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.337"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.345"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05797"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02095"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.142"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.766"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5630"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1781"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.765"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07206"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.69"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5299"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.160"
$ws.Range("E46").Value = "  -5.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.824"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.12"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.058"
$ws.Range("E49").Value = "  -4.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.366"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").Value = "  +0.57%  "
